# Insert a new row at position 403, pushing all existing rows 403..495 down to 404..496.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(403).Insert()

# Populate the freshly inserted row 403 with the new record's data.
$ws.Cells.Item(403, 1).Value = 4
$ws.Cells.Item(403, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(403, 3).Value = "Los Lagos"
$ws.Cells.Item(403, 4).Value = 44754
$ws.Cells.Item(403, 5).Value = 10
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100106
$ws.Cells.Item(403, 8).Value = "Oleaginosos"
$ws.Cells.Item(403, 9).Value = 100106002
$ws.Cells.Item(403, 10).Value = "Palta"
$ws.Cells.Item(403, 11).Value = "Hass"
$ws.Cells.Item(403, 12).Value = "Primera"
$ws.Cells.Item(403, 13).Value = 450
$ws.Cells.Item(403, 14).Value = 27000
$ws.Cells.Item(403, 15).Value = 28000
$ws.Cells.Item(403, 16).Value = 27444
$ws.Cells.Item(403, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(403, 18).Value = "Perú"
$ws.Cells.Item(403, 19).Value = 2744
$ws.Cells.Item(403, 20).Value = 10
